# Edit the "startup" worksheet in the UBC02 stage-of-disease clinical workbook.
# The shared query text in B2 (the "CasesTab" Cypher query) drops its trailing
# "Cohort" projection line, and the now-shorter wrapped cell causes the data
# rows to re-flow at a slightly shorter row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit: remove the trailing `coalesce(co.cohort_description, '') AS `Cohort`` line
# from the CasesTab query stored in B2, keeping everything else identical. ---
$newQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1'] OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value = $newQuery

# --- The wrapped/merged query cells reflow to a shorter row height after the edit. ---
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 290

# --- Scroll the view down one row so row 2 sits at the top of the window. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
